{"js": "// Map of old text -> new text, derived from the OOXML diff.\nconst replacements = [\n  [\"2025-01-02 Thursday\", \"2025-01-03 Friday\"],\n  [\"326\u00f74=\", \"790\u00f74=\"],\n  [\"356\u00f74=\", \"745\u00f74=\"],\n  [\"987\u00f74=\", \"493\u00f74=\"],\n  [\"141\u00f72=\", \"476\u00f76=\"],\n  [\"577\u00f78=\", \"511\u00f79=\"],\n  [\"752\u00f75=\", \"243\u00f72=\"],\n  [\"130\u00f76=\", \"675\u00f79=\"],\n  [\"213\u00f75=\", \"491\u00f78=\"],\n  [\"248\u00f75=\", \"482\u00f78=\"],\n  [\"797\u00f74=\", \"699\u00f78=\"],\n  [\"556\u00f75=\", \"769\u00f79=\"],\n  [\"478\u00f78=\", \"726\u00f74=\"],\n  [\"541\u00f75=\", \"723\u00f75=\"],\n  [\"933\u00f78=\", \"142\u00f75=\"],\n  [\"180\u00f72=\", \"201\u00f79=\"],\n  [\"380\u00f76=\", \"413\u00f78=\"],\n  [\"249\u00f72=\", \"202\u00f72=\"],\n  [\"112\u00f72=\", \"471\u00f78=\"],\n  [\"641\u00f77=\", \"871\u00f78=\"],\n  [\"819\u00f74=\", \"288\u00f78=\"],\n  [\"229\u00f79=\", \"110\u00f75=\"],\n  [\"695\u00f72=\", \"542\u00f73=\"],\n  [\"782\u00f72=\", \"241\u00f73=\"],\n  [\"154\u00f77=\", \"481\u00f73=\"],\n  [\"693\u00f74=\", \"387\u00f72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-01-02 Thursday\", \"2025-01-03 Friday\"),\n    @(\"326\u00f74=\", \"790\u00f74=\"),\n    @(\"356\u00f74=\", \"745\u00f74=\"),\n    @(\"987\u00f74=\", \"493\u00f74=\"),\n    @(\"141\u00f72=\", \"476\u00f76=\"),\n    @(\"577\u00f78=\", \"511\u00f79=\"),\n    @(\"752\u00f75=\", \"243\u00f72=\"),\n    @(\"130\u00f76=\", \"675\u00f79=\"),\n    @(\"213\u00f75=\", \"491\u00f78=\"),\n    @(\"248\u00f75=\", \"482\u00f78=\"),\n    @(\"797\u00f74=\", \"699\u00f78=\"),\n    @(\"556\u00f75=\", \"769\u00f79=\"),\n    @(\"478\u00f78=\", \"726\u00f74=\"),\n    @(\"541\u00f75=\", \"723\u00f75=\"),\n    @(\"933\u00f78=\", \"142\u00f75=\"),\n    @(\"180\u00f72=\", \"201\u00f79=\"),\n    @(\"380\u00f76=\", \"413\u00f78=\"),\n    @(\"249\u00f72=\", \"202\u00f72=\"),\n    @(\"112\u00f72=\", \"471\u00f78=\"),\n    @(\"641\u00f77=\", \"871\u00f78=\"),\n    @(\"819\u00f74=\", \"288\u00f78=\"),\n    @(\"229\u00f79=\", \"110\u00f75=\"),\n    @(\"695\u00f72=\", \"542\u00f73=\"),\n    @(\"782\u00f72=\", \"241\u00f73=\"),\n    @(\"154\u00f77=\", \"481\u00f73=\"),\n    @(\"693\u00f74=\", \"387\u00f72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
